$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 872) into the
# new rows so column A keeps the bold/bordered/centered style (s="1").
$ws.Range("A872").Copy()
$ws.Range("A873:A886").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(873, 1).Value = 871
$ws.Cells.Item(873, 2).Value = 'carl edwards'
$ws.Cells.Item(873, 3).Value = 2.76
$ws.Cells.Item(874, 1).Value = 872
$ws.Cells.Item(874, 2).Value = 'daniel lynch'
$ws.Cells.Item(874, 3).Value = 5.13
$ws.Cells.Item(875, 1).Value = 873
$ws.Cells.Item(875, 2).Value = 'duane underwood'
$ws.Cells.Item(875, 3).Value = 4.4
$ws.Cells.Item(876, 1).Value = 874
$ws.Cells.Item(876, 2).Value = 'frank german'
$ws.Cells.Item(876, 3).Value = 18
$ws.Cells.Item(877, 1).Value = 875
$ws.Cells.Item(877, 2).Value = 'j.t. chargois'
$ws.Cells.Item(877, 3).Value = 2.42
$ws.Cells.Item(878, 1).Value = 876
$ws.Cells.Item(878, 2).Value = 'jaime barría'
$ws.Cells.Item(878, 3).Value = 2.61
$ws.Cells.Item(879, 1).Value = 877
$ws.Cells.Item(879, 2).Value = 'lance mccullers'
$ws.Cells.Item(879, 3).Value = 2.27
$ws.Cells.Item(880, 1).Value = 878
$ws.Cells.Item(880, 2).Value = 'mark leiter'
$ws.Cells.Item(880, 3).Value = 3.99
$ws.Cells.Item(881, 1).Value = 879
$ws.Cells.Item(881, 2).Value = 'matt boyd'
$ws.Cells.Item(881, 3).Value = 1.35
$ws.Cells.Item(882, 1).Value = 880
$ws.Cells.Item(882, 2).Value = 'mike king'
$ws.Cells.Item(882, 3).Value = 2.29
$ws.Cells.Item(883, 1).Value = 881
$ws.Cells.Item(883, 2).Value = 'néstor cortés'
$ws.Cells.Item(883, 3).Value = 2.44
$ws.Cells.Item(884, 1).Value = 882
$ws.Cells.Item(884, 2).Value = 'nick martínez'
$ws.Cells.Item(884, 3).Value = 3.47
$ws.Cells.Item(885, 1).Value = 883
$ws.Cells.Item(885, 2).Value = 'travis lakins'
$ws.Cells.Item(885, 3).Value = 9.58
$ws.Cells.Item(886, 1).Value = 884
$ws.Cells.Item(886, 2).Value = 'vladimir gutiérrez'
$ws.Cells.Item(886, 3).Value = 7.61

Write-Host "Rows added. Dimension should now be A1:C886."
